# Auto-generated script to update Leve profit calculation values
# per scheduled market-data refresh (see commit message).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2681.7273
$ws.Range("I113").Value = 2496.5
$ws.Range("J113").Value = 2787.5715
$ws.Range("K113").Value = 2496.5
$ws.Range("L113").Value = 2787.5715
$ws.Range("M113").Value = 757.5
$ws.Range("N113").Value = -9295.5715

$ws.Range("H138").Value = 2286.535
$ws.Range("I138").Value = 1771.0667
$ws.Range("J138").Value = 2562.6785
$ws.Range("K138").Value = 5313.2001
$ws.Range("L138").Value = 7688.0355
$ws.Range("M138").Value = -173.2001
$ws.Range("N138").Value = -17968.0355

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13123721
$ws.Range("I32").Value = 14761709
$ws.Range("K32").Value = 14761709
$ws.Range("M32").Value = -14761422

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1179.95
$ws.Range("I80").Value = 1889.7273
$ws.Range("K80").Value = 1889.7273
$ws.Range("M80").Value = -891.7273

$ws.Range("H83").Value = 1179.95
$ws.Range("I83").Value = 1889.7273
$ws.Range("K83").Value = 9448.636500000001
$ws.Range("M83").Value = -4456.636500000001

$ws.Range("H94").Value = 1649.9166
$ws.Range("I94").Value = 1484.1428
$ws.Range("K94").Value = 1484.1428
$ws.Range("M94").Value = -1033.1428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4570.216
$ws.Range("I31").Value = 1163.5
$ws.Range("K31").Value = 1163.5
$ws.Range("M31").Value = -868.5

$ws.Range("H34").Value = 4570.216
$ws.Range("I34").Value = 1163.5
$ws.Range("K34").Value = 1163.5
$ws.Range("M34").Value = -961.5

$ws.Range("H51").Value = 16199.4
$ws.Range("J51").Value = 16199.4
$ws.Range("L51").Value = 16199.4
$ws.Range("N51").Value = -17671.4

$ws.Range("H58").Value = 1059.1277
$ws.Range("I58").Value = 778.4194
$ws.Range("K58").Value = 778.4194
$ws.Range("M58").Value = -575.4194

$ws.Range("H61").Value = 16199.4
$ws.Range("J61").Value = 16199.4
$ws.Range("L61").Value = 16199.4
$ws.Range("N61").Value = -16895.4

$ws.Range("H99").Value = 2108.6956
$ws.Range("I99").Value = 1901
$ws.Range("J99").Value = 2118.1365
$ws.Range("K99").Value = 1901
$ws.Range("L99").Value = 2118.1365
$ws.Range("M99").Value = -403
$ws.Range("N99").Value = -5114.136500000001

$ws.Range("H105").Value = 414.66666
$ws.Range("I105").Value = 400.4
$ws.Range("J105").Value = 700
$ws.Range("K105").Value = 400.4
$ws.Range("L105").Value = 700
$ws.Range("M105").Value = 1346.6
$ws.Range("N105").Value = -4194

$ws.Range("H122").Value = 1981.1
$ws.Range("I122").Value = 1874
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 5622
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -3172
$ws.Range("N122").Value = -10900

$ws.Range("H126").Value = 2108.6956
$ws.Range("I126").Value = 1901
$ws.Range("J126").Value = 2118.1365
$ws.Range("K126").Value = 5703
$ws.Range("L126").Value = 6354.4095
$ws.Range("M126").Value = -3233
$ws.Range("N126").Value = -11294.4095

$ws.Range("H134").Value = 2770.1455
$ws.Range("I134").Value = 2901.16
$ws.Range("J134").Value = 1460
$ws.Range("K134").Value = 8703.48
$ws.Range("L134").Value = 4380
$ws.Range("M134").Value = -6168.48
$ws.Range("N134").Value = -9450

$ws.Range("H136").Value = 1059.1277
$ws.Range("I136").Value = 778.4194
$ws.Range("K136").Value = 2335.2582
$ws.Range("M136").Value = 214.7417999999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 1871.875
$ws.Range("I58").Value = 475
$ws.Range("J58").Value = 2071.4285
$ws.Range("K58").Value = 1425
$ws.Range("L58").Value = 6214.2855
$ws.Range("M58").Value = -1297
$ws.Range("N58").Value = -6470.2855

$ws.Range("H80").Value = 3338.2727
$ws.Range("I80").Value = 1601
$ws.Range("J80").Value = 3724.3333
$ws.Range("K80").Value = 4803
$ws.Range("L80").Value = 11172.9999
$ws.Range("M80").Value = -3867
$ws.Range("N80").Value = -13044.9999

$ws.Range("H83").Value = 3338.2727
$ws.Range("I83").Value = 1601
$ws.Range("J83").Value = 3724.3333
$ws.Range("K83").Value = 14409
$ws.Range("L83").Value = 33518.9997
$ws.Range("M83").Value = -9729
$ws.Range("N83").Value = -42878.9997

$ws.Range("H139").Value = 3356.6052
$ws.Range("I139").Value = 1697.6111
$ws.Range("J139").Value = 4849.7
$ws.Range("K139").Value = 5092.8333
$ws.Range("L139").Value = 14549.1
$ws.Range("M139").Value = 47.16669999999976
$ws.Range("N139").Value = -24829.1

$ws.Range("H140").Value = 1544.9736
$ws.Range("I140").Value = 1040.36
$ws.Range("J140").Value = 2515.3845
$ws.Range("K140").Value = 3121.08
$ws.Range("L140").Value = 7546.1535
$ws.Range("M140").Value = 2058.92
$ws.Range("N140").Value = -17906.1535

$ws.Range("H141").Value = 5651.7
$ws.Range("I141").Value = 3941.077
$ws.Range("J141").Value = 8828.571
$ws.Range("K141").Value = 11823.231
$ws.Range("L141").Value = 26485.713
$ws.Range("M141").Value = -6643.231
$ws.Range("N141").Value = -36845.713

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 14287014
$ws.Range("I16").Value = 1600.2858
$ws.Range("J16").Value = 47619644
$ws.Range("K16").Value = 1600.2858
$ws.Range("L16").Value = 47619644
$ws.Range("M16").Value = -1430.2858
$ws.Range("N16").Value = -47619984

$ws.Range("H132").Value = 2366.4375
$ws.Range("I132").Value = 2161.0186
$ws.Range("K132").Value = 6483.0558
$ws.Range("M132").Value = -3953.0558

$ws.Range("H136").Value = 11113397
$ws.Range("I136").Value = 3016.8333
$ws.Range("J136").Value = 18520316
$ws.Range("K136").Value = 9050.499899999999
$ws.Range("L136").Value = 55560948
$ws.Range("M136").Value = -6500.499899999999
$ws.Range("N136").Value = -55566048

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4066306
$ws.Range("I132").Value = 1219.7037
$ws.Range("J132").Value = 11906115
$ws.Range("K132").Value = 3659.1111
$ws.Range("L132").Value = 35718345
$ws.Range("M132").Value = -1129.1111
$ws.Range("N132").Value = -35723405

$ws.Range("H136").Value = 2576.7334
$ws.Range("I136").Value = 2168.6924
$ws.Range("J136").Value = 5229
$ws.Range("K136").Value = 6506.0772
$ws.Range("L136").Value = 15687
$ws.Range("M136").Value = -3956.0772
$ws.Range("N136").Value = -20787
